$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.462.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.24%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.022.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.25%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'254.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.07%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'56.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -6.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.56%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.85%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.101"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.25%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.320.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.24%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.815"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.78%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'21.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.40%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.017.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.91%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.334.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.27%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.81%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.98%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.73%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'228.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.76%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'163.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.63%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.34%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "'  -10.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.94%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.83%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0664"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +7.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.62%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.52%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.22%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.64%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.57%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.21%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.50%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.22%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.396.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.13%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'15.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.39%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'90.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.70%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.96%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.212.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.26%  "
$ws.Range("E51").Style = "Normal"
